$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "311.91"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.81%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.69"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.52%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.125"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.92%"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.81%"

$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.411"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.70%"

$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.911"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.23%"

$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "8.271"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.11%"

$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.838"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-7.51%"

$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9184"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.22%"

$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1170"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-4.17%"

$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1927"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "3.19%"

$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09021"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.37%"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03313"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-2.52%"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09598"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.10%"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001375"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.40%"

$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006057"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.26%"

$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.535"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.34%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3442"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.99%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "4.95%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1285"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.29%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2587"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "3.94%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04356"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.72%"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.02%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004669"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "10.37%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.68%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003987"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02258"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.66%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05079"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.78%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007461"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.68%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009040"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-8.76%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1354"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.06%"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-5.38%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008624"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-12.20%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006573"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.62%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000749"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0009993"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-23.17%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003005"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.34%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002098"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001998"
